$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '68.141.33'
$ws.Cells.Item(2, 5).Value2 = '  +0.43%  '

Set-TextValue 3 4 '3.677.41'
$ws.Cells.Item(3, 5).Value2 = '  -0.52%  '

Set-TextValue 4 4 '0.999'
$ws.Cells.Item(4, 5).Value2 = '  -0.17%  '

Set-TextValue 5 4 '601.42'
$ws.Cells.Item(5, 5).Value2 = '  +4.20%  '

Set-TextValue 6 4 '191.71'
$ws.Cells.Item(6, 5).Value2 = '  +9.70%  '

$ws.Cells.Item(7, 5).Value2 = '  +0.17%  '

Set-TextValue 8 4 '1.00'
$ws.Cells.Item(8, 5).Value2 = '  +0.23%  '

Set-TextValue 9 4 '0.708'
$ws.Cells.Item(9, 5).Value2 = '  +0.83%  '

Set-TextValue 10 4 '58.26'
$ws.Cells.Item(10, 5).Value2 = '  +13.49%  '

Set-TextValue 11 4 '0.154'
$ws.Cells.Item(11, 5).Value2 = '  -4.60%  '

Set-TextValue 12 4 '0.0000277'
$ws.Cells.Item(12, 5).Value2 = '  -3.95%  '

Set-TextValue 13 4 '10.26'
$ws.Cells.Item(13, 5).Value2 = '  -1.02%  '

Set-TextValue 14 4 '4.255.97'
$ws.Cells.Item(14, 5).Value2 = '  -0.72%  '

Set-TextValue 15 4 '3.672.23'
$ws.Cells.Item(15, 5).Value2 = '  -0.60%  '

$ws.Cells.Item(16, 5).Value2 = '  +0.91%  '

Set-TextValue 17 4 '19.07'
$ws.Cells.Item(17, 5).Value2 = '  -1.01%  '

Set-TextValue 18 4 '1.13'
$ws.Cells.Item(18, 5).Value2 = '  +1.28%  '

Set-TextValue 19 4 '67.898.04'
$ws.Cells.Item(19, 5).Value2 = '  +0.37%  '

Set-TextValue 20 4 '12.57'
$ws.Cells.Item(20, 5).Value2 = '  -1.43%  '

Set-TextValue 21 4 '402.44'
$ws.Cells.Item(21, 5).Value2 = '  -0.43%  '

$ws.Cells.Item(22, 5).Value2 = '  +0.51%  '

Set-TextValue 23 4 '88.47'
$ws.Cells.Item(23, 5).Value2 = '  +0.69%  '

Set-TextValue 24 4 '11.38'
$ws.Cells.Item(24, 5).Value2 = '  +6.34%  '

Set-TextValue 25 4 '2.98'
$ws.Cells.Item(25, 5).Value2 = '  -1.21%  '

Set-TextValue 26 4 '12.65'
$ws.Cells.Item(26, 5).Value2 = '  +0.02%  '

$ws.Cells.Item(27, 5).Value2 = '  -0.04%  '

Set-TextValue 28 4 '3.69'
$ws.Cells.Item(28, 5).Value2 = '  -1.70%  '

Set-TextValue 29 4 '9.38'
$ws.Cells.Item(29, 5).Value2 = '  -0.54%  '

Set-TextValue 30 4 '32.12'
$ws.Cells.Item(30, 5).Value2 = '  -0.82%  '

Set-TextValue 31 4 '7.61'
$ws.Cells.Item(31, 5).Value2 = '  +3.07%  '

Set-TextValue 32 4 '45.86'
$ws.Cells.Item(32, 5).Value2 = '  +7.27%  '

Set-TextValue 33 4 '12.44'
$ws.Cells.Item(33, 5).Value2 = '  +0.66%  '

Set-TextValue 34 4 '67.24'
$ws.Cells.Item(34, 5).Value2 = '  +3.76%  '

$ws.Cells.Item(35, 2).Value2 = 'Hedera'
$ws.Cells.Item(35, 3).Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 35 4 '0.118'
$ws.Cells.Item(35, 5).Value2 = '  +2.55%  '

$ws.Cells.Item(36, 2).Value2 = 'Bittensor'
$ws.Cells.Item(36, 3).Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 36 4 '617.84'
$ws.Cells.Item(36, 5).Value2 = '  +1.32%  '

Set-TextValue 37 4 '1.00'
$ws.Cells.Item(37, 5).Value2 = '  -0.09%  '

Set-TextValue 38 4 '0.400'
$ws.Cells.Item(38, 5).Value2 = '  +1.69%  '

$ws.Cells.Item(39, 2).Value2 = 'PEPE'
$ws.Cells.Item(39, 3).Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 39 4 '0.0₃0782'
$ws.Cells.Item(39, 5).Value2 = '  -11.08%  '

$ws.Cells.Item(40, 2).Value2 = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 40 4 '0.999'
$ws.Cells.Item(40, 5).Value2 = '  -0.09%  '

Set-TextValue 41 4 '0.136'
$ws.Cells.Item(41, 5).Value2 = '  +0.94%  '

Set-TextValue 42 4 '2.93'
$ws.Cells.Item(42, 5).Value2 = '  -1.34%  '

Set-TextValue 43 4 '0.0430'
$ws.Cells.Item(43, 5).Value2 = '  -0.25%  '

Set-TextValue 44 4 '2.58'
$ws.Cells.Item(44, 5).Value2 = '  -6.06%  '

Set-TextValue 45 4 '2.838.85'
$ws.Cells.Item(45, 5).Value2 = '  +1.74%  '

$ws.Cells.Item(46, 5).Value2 = '  +2.83%  '

$ws.Cells.Item(47, 5).Value2 = '  +6.05%  '

Set-TextValue 48 4 '9.00'
$ws.Cells.Item(48, 5).Value2 = '  -1.81%  '

Set-TextValue 49 4 '144.94'
$ws.Cells.Item(49, 5).Value2 = '  +4.77%  '

Set-TextValue 50 4 '2.65'
$ws.Cells.Item(50, 5).Value2 = '  -1.29%  '

Set-TextValue 51 4 '2.52'
$ws.Cells.Item(51, 5).Value2 = '  -10.94%  '
